$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the D2:E51 range (skipping header) so numeric-looking
# strings like "0.9999" are stored as text, matching the source data.
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.778.24'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '1.796.54'
$ws.Range("E3").Value = '  -1.29%  '
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '309.13'
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").Value = '0.4402'
$ws.Range("E7").Value = '  +4.38%  '
$ws.Range("D8").Value = '0.3667'
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("D9").Value = '0.07325'
$ws.Range("E9").Value = '  +1.63%  '
$ws.Range("D10").Value = '0.8536'
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("E11").Value = '  -1.05%  '
$ws.Range("D12").Value = '1.926.63'
$ws.Range("E12").Value = '  +5.89%  '
$ws.Range("D13").Value = '6.620'
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("D14").Value = '92.05'
$ws.Range("E14").Value = '  +2.65%  '
$ws.Range("D15").Value = '0.07037'
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = '5.261'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = '0.000008640'
$ws.Range("E18").Value = '  -1.55%  '
$ws.Range("D19").Value = '0.9998'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("E20").Value = '  -1.19%  '
$ws.Range("D21").Value = '26.807.61'
$ws.Range("E21").Value = '  -1.07%  '
$ws.Range("D22").Value = '5.134'
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("D23").Value = '10.78'
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("D24").Value = '1.974'
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = '151.42'
$ws.Range("E25").Value = '  -0.23%  '
$ws.Range("D26").Value = '2.196'
$ws.Range("E26").Value = '  -1.12%  '
$ws.Range("D27").Value = '18.37'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("D28").Value = '5.175'
$ws.Range("E28").Value = '  -1.21%  '
$ws.Range("D29").Value = '116.71'
$ws.Range("E29").Value = '  +0.46%  '
$ws.Range("D30").Value = '0.08769'
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("D31").Value = '0.7357'
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("D32").Value = '1.150'
$ws.Range("E32").Value = '  -2.09%  '
$ws.Range("E33").Value = '  -0.43%  '
$ws.Range("D34").Value = '4.419'
$ws.Range("E34").Value = '  +0.16%  '
$ws.Range("D35").Value = '0.9995'
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").Value = '1.083'
$ws.Range("E36").Value = '  -0.65%  '
$ws.Range("D37").Value = '0.01952'
$ws.Range("E37").Value = '  +0.19%  '
$ws.Range("D39").Value = '0.5201'
$ws.Range("E39").Value = '  +3.44%  '
$ws.Range("E40").Value = '  -4.55%  '
$ws.Range("D41").Value = '2.803'
$ws.Range("E41").Value = '  -2.49%  '
$ws.Range("D42").Value = '0.1672'
$ws.Range("E42").Value = '  -0.80%  '
$ws.Range("D43").Value = '0.4990'
$ws.Range("E43").Value = '  +6.10%  '
$ws.Range("D44").Value = '8.404'
$ws.Range("E44").Value = '  -1.89%  '
$ws.Range("D45").Value = '1.964'
$ws.Range("E45").Value = '  +4.37%  '
$ws.Range("D46").Value = '10.29'
$ws.Range("E46").Value = '  -2.45%  '
$ws.Range("D47").Value = '104.87'
$ws.Range("E47").Value = '  -1.18%  '
$ws.Range("D48").Value = '0.9991'
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").Value = '1.655'
$ws.Range("E49").Value = '  +0.63%  '
$ws.Range("E50").Value = '  -0.70%  '
$ws.Range("D51").Value = '0.9132'
$ws.Range("E51").Value = '  +1.40%  '

# Restore the default (unstyled) cell style now that the text values are set,
# so the cells keep no explicit style index, same as before the edit.
$fmtRange.Style = "Normal"

